$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.155.69'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '2.571.68'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.72'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.21'
$ws.Range("E6").Value = '  -2.35%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.61'
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.13'
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").Value = '3.035.62'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '63.111.19'
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000145'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").Value = '2.580.31'
$ws.Range("E17").Value = '  +1.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.04'
$ws.Range("E18").Value = '  -2.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '340.68'
$ws.Range("E19").Value = '  -0.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.33'
$ws.Range("E20").Value = '  -1.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.63'
$ws.Range("E21").Value = '  -3.43%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.73'
$ws.Range("E23").Value = '  +3.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.90'
$ws.Range("E24").Value = '  +1.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.55'
$ws.Range("E25").Value = '  +4.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.60'
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.164'
$ws.Range("E27").Value = '  -3.25%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.92'
$ws.Range("E29").Value = '  -2.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.19'
$ws.Range("E30").Value = '  -2.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.93'
$ws.Range("E31").Value = '  -2.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '470.11'
$ws.Range("E32").Value = '  +2.09%  '
$ws.Range("D33").Value = '0.0₃0798'
$ws.Range("E33").Value = '  -3.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.67'
$ws.Range("E34").Value = '  +3.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '176.04'
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("E37").Value = '  -1.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.84'
$ws.Range("E38").Value = '  -1.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.53'
$ws.Range("E39").Value = '  +0.20%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("E41").Value = '  -3.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.17'
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '157.58'
$ws.Range("E43").Value = '  +4.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.68'
$ws.Range("E44").Value = '  -3.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.20'
$ws.Range("E45").Value = '  +1.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.631'
$ws.Range("E46").Value = '  +2.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0538'
$ws.Range("E47").Value = '  -1.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0962'
$ws.Range("E48").Value = '  -1.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0236'
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.11'
$ws.Range("E50").Value = '  -1.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.37'
$ws.Range("E51").Value = '  -0.12%  '
